$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Claims")

# Step 1: move the "less than 0" highlight rule from column E (Sent) to column D (Bill), reusing dxfId=0
$fc = $ws.Range("E1:E1048576").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D1:D1048576"))

# Step 2: remove the old "Sent" column entirely -- "Paid" (old F) slides left into E
$ws.Range("E1").EntireColumn.Delete()

# Step 3: rename headers C->Range, D->Bill, and give D the currency format/boldness used by Paid (E)
$ws.Range("C1").Value = "Range"
$ws.Range("D1").Value = "Bill"
$ws.Range("D1").NumberFormat = $ws.Range("E1").NumberFormat
$ws.Range("D1").Font.Bold = $ws.Range("E1").Font.Bold

# Step 4: update existing claim rows with the new consolidated Range + Bill layout
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Value = "5/1/24 - 5/17/24"
$ws.Range("D2").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("D2").Value = 1300
$ws.Range("E2").Clear()
$ws.Range("B2").Value = "McGee, Test"
$ws.Range("A2").Value = 45429

$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = "5/1/24 - 5/17/24"
$ws.Range("D3").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("D3").Value = 400
$ws.Range("E3").Clear()
$ws.Range("B3").Value = "Anna, Mary"
$ws.Range("A3").Value = 45429

# Step 5: pre-format the next two rows' date column the same way the sheet already did for row 4
$ws.Range("C4").Clear()
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Output "done"
